$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to remain plain text (avoid Excel
# auto-converting numeric-looking strings like "236.18" or
# "0.000006622" into real numbers / scientific notation).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.016.64'
$ws.Range("D3").Value = '1.646.97'
$ws.Range("D4").Value = '0.9994'
$ws.Range("D5").Value = '236.18'
$ws.Range("D7").Value = '0.4839'
$ws.Range("D8").Value = '0.2602'
$ws.Range("D9").Value = '0.06006'
$ws.Range("D10").Value = '0.07194'
$ws.Range("D11").Value = '1.645.23'
$ws.Range("D12").Value = '14.78'
$ws.Range("D13").Value = '0.6209'
$ws.Range("D14").Value = '4.519'
$ws.Range("D15").Value = '72.81'
$ws.Range("D16").Value = '1.001'
$ws.Range("D17").Value = '0.9991'
$ws.Range("D18").Value = '25.002.91'
$ws.Range("D19").Value = '11.42'
$ws.Range("D20").Value = '0.000006622'
$ws.Range("D21").Value = '4.519'
$ws.Range("D22").Value = '1.856.65'
$ws.Range("D23").Value = '8.622'
$ws.Range("D24").Value = '5.289'
$ws.Range("D25").Value = '132.02'
$ws.Range("D26").Value = '14.89'
$ws.Range("D27").Value = '1.401'
$ws.Range("D28").Value = '102.86'
$ws.Range("D29").Value = '1.672'
$ws.Range("D30").Value = '3.756'
$ws.Range("D31").Value = '0.07858'
$ws.Range("D33").Value = '0.04498'
$ws.Range("D34").Value = '0.9994'
$ws.Range("D35").Value = '2.595'
$ws.Range("D36").Value = '0.9339'
$ws.Range("D37").Value = '0.5821'
$ws.Range("D38").Value = '2.576'
$ws.Range("D39").Value = '0.01566'
$ws.Range("D40").Value = '0.8488'
$ws.Range("D41").Value = '0.9995'
$ws.Range("D43").Value = '98.18'
$ws.Range("D44").Value = '0.3720'
$ws.Range("D45").Value = '4.783'
$ws.Range("D46").Value = '0.1150'
$ws.Range("D47").Value = '6.121'
$ws.Range("D48").Value = '0.05189'
$ws.Range("D49").Value = '29.77'
$ws.Range("D50").Value = '1.000'
$ws.Range("D51").Value = '50.41'

# Coin / Link / Volume columns are unambiguous text already
# (letters, URLs, or percent-signs with surrounding spaces), so
# no NumberFormat coercion is required for these.
$ws.Range("E2").Value = '  -3.71%  '
$ws.Range("E3").Value = '  -5.45%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("E5").Value = '  -5.72%  '
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("E7").Value = '  -6.14%  '
$ws.Range("E8").Value = '  -5.59%  '
$ws.Range("E9").Value = '  -2.98%  '
$ws.Range("E10").Value = '  -0.46%  '
$ws.Range("E11").Value = '  -5.60%  '
$ws.Range("E12").Value = '  -2.40%  '
$ws.Range("E13").Value = '  -4.56%  '
$ws.Range("E14").Value = '  -2.54%  '
$ws.Range("E15").Value = '  -6.21%  '
$ws.Range("E16").Value = '  +0.03%  '
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  -3.83%  '
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("E20").Value = '  -2.81%  '
$ws.Range("E21").Value = '  +5.77%  '
$ws.Range("E22").Value = '  -5.56%  '
$ws.Range("E23").Value = '  -0.76%  '
$ws.Range("E24").Value = '  -1.59%  '
$ws.Range("E25").Value = '  -3.07%  '
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("E27").Value = '  -7.65%  '
$ws.Range("E28").Value = '  -2.82%  '
$ws.Range("E29").Value = '  -6.16%  '
$ws.Range("E30").Value = '  -5.28%  '
$ws.Range("E31").Value = '  -4.48%  '
$ws.Range("E32").Value = '  -2.01%  '
$ws.Range("E33").Value = '  -4.40%  '
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("E35").Value = '  -2.38%  '
$ws.Range("E36").Value = '  -6.27%  '
$ws.Range("E37").Value = '  -6.72%  '
$ws.Range("E38").Value = '  -5.62%  '
$ws.Range("E39").Value = '  -3.23%  '
$ws.Range("E40").Value = '  +11.54%  '
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("E43").Value = '  -2.31%  '
$ws.Range("E44").Value = '  -3.39%  '
$ws.Range("E45").Value = '  -4.98%  '
$ws.Range("E46").Value = '  +1.66%  '
$ws.Range("E47").Value = '  -3.09%  '
$ws.Range("E48").Value = '  -0.77%  '
$ws.Range("E49").Value = '  -3.30%  '
$ws.Range("B50").Value = 'TrueUSD'
$ws.Range("C50").Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E51").Value = '  -9.41%  '
